$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price number format to show 4 decimal places instead of 2
# (restricted to the populated price cells so we don't materialize blank
# styled cells across the whole column)
$fmt = "[$$-409]#,##0.0000;[RED]\-[$$-409]#,##0.0000"
$ws.Range("H1:J89").NumberFormat = $fmt
$ws.Range("K1").NumberFormat = $fmt
$ws.Range("K85:K89").NumberFormat = $fmt

# Update J60 value from 367 to 0.367
$ws.Cells.Item(60, 10).Value = 0.367

# Delete row 90 (ATLAS_SENSOR / F91 row) and shift rows below up.
$ws.Rows.Item(90).Delete()

Write-Host "Done"
